$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.213.73"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "2.307.21"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.96%  "
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.42"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.23%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.118"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +14.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("D15").Value = "2.665.14"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "2.286.01"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.817"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.91%  "
$ws.Range("D18").Value = "43.097.11"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.58"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +11.46%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.45"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("E24").Value = "  +15.49%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.08"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.25"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("E34").Value = "  +4.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.55"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.09"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0692"
$ws.Range("D37").ClearFormats()
$ws.Range("E38").Value = "  +3.92%  "
$ws.Range("E39").Value = "  +5.70%  "
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  -3.22%  "
$ws.Range("D43").Value = "1.993.35"
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.11"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.72"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.86"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.01%  "
$ws.Range("D49").Value = "2.533.20"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.57"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.68%  "
